$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.294.46"
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "2.018.13"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'249.80"
$ws.Range("E5").Value = "  +2.55%  "
$ws.Range("D6").Value = "'0.643"
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("D7").Value = "'63.29"
$ws.Range("E7").Value = "  +20.08%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'59.40"
$ws.Range("E9").Value = "  -3.13%  "
$ws.Range("D10").Value = "'0.370"
$ws.Range("E10").Value = "  +3.49%  "
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("E12").Value = "  -0.42%  "
$ws.Range("D13").Value = "'0.947"
$ws.Range("E13").Value = "  +0.83%  "
$ws.Range("D14").Value = "'15.05"
$ws.Range("E14").Value = "  +5.68%  "
$ws.Range("D15").Value = "2.311.88"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "'5.44"
$ws.Range("E16").Value = "  +3.25%  "
$ws.Range("D17").Value = "'19.48"
$ws.Range("E17").Value = "  +17.40%  "
$ws.Range("D18").Value = "2.023.83"
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").Value = "36.207.14"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "'72.20"
$ws.Range("E20").Value = "  +2.10%  "
$ws.Range("D21").Value = "0.0₃0858"
$ws.Range("E21").Value = "  +1.87%  "
$ws.Range("D22").Value = "'5.31"
$ws.Range("E22").Value = "  +3.98%  "
$ws.Range("D23").Value = "'234.62"
$ws.Range("E23").Value = "  -0.57%  "
$ws.Range("D24").Value = "'2.68"
$ws.Range("E24").Value = "  +22.11%  "
$ws.Range("E25").Value = "  +0.07%  "
$ws.Range("D26").Value = "'2.29"
$ws.Range("E26").Value = "  -2.57%  "
$ws.Range("D27").Value = "'9.65"
$ws.Range("E27").Value = "  +7.08%  "
$ws.Range("D28").Value = "'166.13"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").Value = "'19.66"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "'0.121"
$ws.Range("E30").Value = "  +1.30%  "
$ws.Range("D31").Value = "'5.17"
$ws.Range("E31").Value = "  +6.42%  "
$ws.Range("D32").Value = "'1.19"
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "'0.106"
$ws.Range("E33").Value = "  +25.13%  "
$ws.Range("D34").Value = "'0.0609"
$ws.Range("E34").Value = "  +3.73%  "
$ws.Range("D35").Value = "'4.50"
$ws.Range("E35").Value = "  +4.11%  "
$ws.Range("D36").Value = "'2.48"
$ws.Range("E36").Value = "  +14.14%  "
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "'5.81"
$ws.Range("E39").Value = "  +20.54%  "
$ws.Range("B40").Value = "Cronos"
$ws.Range("C40").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D40").Value = "'0.106"
$ws.Range("E40").Value = "  +19.99%  "
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.22"
$ws.Range("E41").Value = "  +0.62%  "
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.90"
$ws.Range("E42").Value = "  +1.66%  "
$ws.Range("D43").Value = "'0.0217"
$ws.Range("E43").Value = "  +2.92%  "
$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'17.05"
$ws.Range("E44").Value = "  +9.63%  "
$ws.Range("D45").Value = "'1.13"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("D46").Value = "'94.84"
$ws.Range("E46").Value = "  +3.07%  "
$ws.Range("D47").Value = "'7.80"
$ws.Range("E47").Value = "  +5.19%  "
$ws.Range("D48").Value = "1.386.20"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "'2.36"
$ws.Range("E50").Value = "  +5.89%  "
$ws.Range("D51").Value = "'47.23"
$ws.Range("E51").Value = "  +5.52%  "
